$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 154; $r++) {
    $cell = $ws.Cells.Item($r, 5)  # column E = 5
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $cell.Value = -$val
    }
}
